$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-04-11 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-04-12 Friday", 2)

$d.Content.Find.Execute("12÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "69÷9=", 2)
$d.Content.Find.Execute("45÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "82÷9=", 2)
$d.Content.Find.Execute("79÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "75÷6=", 2)
$d.Content.Find.Execute("61÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "60÷5=", 2)
$d.Content.Find.Execute("74÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "82÷6=", 2)
$d.Content.Find.Execute("52÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "17÷6=", 2)
$d.Content.Find.Execute("51÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "82÷4=", 2)
$d.Content.Find.Execute("40÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "46÷4=", 2)
$d.Content.Find.Execute("50÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "58÷4=", 2)
$d.Content.Find.Execute("10÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "14÷4=", 2)
$d.Content.Find.Execute("63÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "90÷4=", 2)
$d.Content.Find.Execute("69÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "21÷6=", 2)
$d.Content.Find.Execute("30÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "99÷6=", 2)
$d.Content.Find.Execute("53÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "56÷8=", 2)
$d.Content.Find.Execute("32÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "83÷9=", 2)
$d.Content.Find.Execute("11÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "87÷8=", 2)
$d.Content.Find.Execute("76÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "44÷9=", 2)
$d.Content.Find.Execute("59÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "31÷8=", 2)
$d.Content.Find.Execute("51÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "89÷4=", 2)
$d.Content.Find.Execute("88÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "22÷7=", 2)
$d.Content.Find.Execute("78÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "97÷3=", 2)
$d.Content.Find.Execute("63÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "52÷6=", 2)
$d.Content.Find.Execute("93÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "13÷8=", 2)
$d.Content.Find.Execute("80÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "71÷7=", 2)
$d.Content.Find.Execute("94÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "88÷6=", 2)
